$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new time-record row (row 12) mirroring the existing rows.
$row = 12

$ws.Cells.Item($row, 1).Value = "30.11.2019"
$ws.Cells.Item($row, 2).Value = 0.5
$ws.Cells.Item($row, 3).Value = 0.52083333333333337
$ws.Cells.Item($row, 4).Formula = "=C12-B12"
$ws.Cells.Item($row, 5).Value = "IO Control Unit"
$ws.Cells.Item($row, 6).Value = "Commenting"

# Match the time number formatting used by the neighbouring rows (6-11).
$ws.Cells.Item($row, 2).NumberFormat = "h:mm"
$ws.Cells.Item($row, 3).NumberFormat = "h:mm"
$ws.Cells.Item($row, 4).NumberFormat = "[$]hh:mm;@"

$ws.Range("A13").Select() | Out-Null
